$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.8181043333333333"
$ws.Range("H2").Value = [double]"2.454313"
$ws.Range("I2").Value = [double]"0.003510462371845219"
$ws.Range("J2").Value = [double]"0.003510462371845219"
$ws.Range("M2").Value = [double]"8.906151666666668"
$ws.Range("N2").Value = [double]"26.718455"
$ws.Range("O2").Value = [double]"0.1245005002255258"
$ws.Range("P2").Value = [double]"0.1245005002255258"
$ws.Range("Q2").Value = [double]"7.28616127182389"
$ws.Range("R2").Value = [double]"65.575451446415"
$ws.Range("S2").Value = [double]"0.0004370543213176155"
$ws.Range("T2").Value = [double]"0.0004370543213176155"

$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.8181043333333333"
$ws.Range("H3").Value = [double]"2.454313"
$ws.Range("I3").Value = [double]"0.003510462371845219"
$ws.Range("J3").Value = [double]"0.003510462371845219"
$ws.Range("M3").Value = [double]"0.4515893333333333"
$ws.Range("N3").Value = [double]"1.354768"
$ws.Range("O3").Value = [double]"0.006312838586270617"
$ws.Range("P3").Value = [double]"0.006312838586270617"
$ws.Range("Q3").Value = [double]"0.3694471904871111"
$ws.Range("R3").Value = [double]"3.325024714384"
$ws.Range("S3").Value = [double]"2.216098231663557E-05"
$ws.Range("T3").Value = [double]"2.216098231663557E-05"

$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.8181043333333333"
$ws.Range("H4").Value = [double]"2.454313"
$ws.Range("I4").Value = [double]"0.003510462371845219"
$ws.Range("J4").Value = [double]"0.003510462371845219"
$ws.Range("M4").Value = [double]"7.781650666666667"
$ws.Range("N4").Value = [double]"23.344952"
$ws.Range("O4").Value = [double]"0.1087809232135948"
$ws.Range("P4").Value = [double]"0.1087809232135948"
$ws.Range("Q4").Value = [double]"6.366202130886222"
$ws.Range("R4").Value = [double]"57.295819177976"
$ws.Range("S4").Value = [double]"0.0003818713377159087"
$ws.Range("T4").Value = [double]"0.0003818713377159087"

$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.8181043333333333"
$ws.Range("H5").Value = [double]"2.454313"
$ws.Range("I5").Value = [double]"0.003510462371845219"
$ws.Range("J5").Value = [double]"0.003510462371845219"
$ws.Range("M5").Value = [double]"3.892567333333333"
$ws.Range("N5").Value = [double]"11.677702"
$ws.Range("O5").Value = [double]"0.05441481329981927"
$ws.Range("P5").Value = [double]"0.05441481329981927"
$ws.Range("Q5").Value = [double]"3.184526203191778"
$ws.Range("R5").Value = [double]"28.660735828726"
$ws.Range("S5").Value = [double]"0.0001910211545599984"
$ws.Range("T5").Value = [double]"0.0001910211545599983"

$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"0.8181043333333333"
$ws.Range("H6").Value = [double]"2.454313"
$ws.Range("I6").Value = [double]"0.003510462371845219"
$ws.Range("J6").Value = [double]"0.003510462371845219"
$ws.Range("M6").Value = [double]"43.49559133333333"
$ws.Range("N6").Value = [double]"130.486774"
$ws.Range("O6").Value = [double]"0.6080317382054886"
$ws.Range("P6").Value = [double]"0.6080317382054886"
$ws.Range("Q6").Value = [double]"35.58393175069578"
$ws.Range("R6").Value = [double]"320.255385756262"
$ws.Range("S6").Value = [double]"0.002134472537858011"
$ws.Range("T6").Value = [double]"0.002134472537858011"

$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"0.8181043333333333"
$ws.Range("H7").Value = [double]"2.454313"
$ws.Range("I7").Value = [double]"0.003510462371845219"
$ws.Range("J7").Value = [double]"0.003510462371845219"
$ws.Range("M7").Value = [double]"7.007517"
$ws.Range("N7").Value = [double]"21.022551"
$ws.Range("O7").Value = [double]"0.09795918646930096"
$ws.Range("P7").Value = [double]"0.09795918646930096"
$ws.Range("Q7").Value = [double]"5.732880023607"
$ws.Range("R7").Value = [double]"51.595920212463"
$ws.Range("S7").Value = [double]"0.0003438820380770504"
$ws.Range("T7").Value = [double]"0.0003438820380770504"

$ws.Range("I8").Value = [double]"0.00984519497906274"
$ws.Range("J8").Value = [double]"0.00984519497906274"
$ws.Range("M8").Value = [double]"8.906151666666668"
$ws.Range("N8").Value = [double]"26.718455"
$ws.Range("O8").Value = [double]"0.1245005002255258"
$ws.Range("P8").Value = [double]"0.1245005002255258"
$ws.Range("Q8").Value = [double]"20.43425360297945"
$ws.Range("R8").Value = [double]"183.908282426815"
$ws.Range("S8").Value = [double]"0.001225731699711146"
$ws.Range("T8").Value = [double]"0.001225731699711146"

$ws.Range("I9").Value = [double]"0.00984519497906274"
$ws.Range("J9").Value = [double]"0.00984519497906274"
$ws.Range("M9").Value = [double]"0.4515893333333333"
$ws.Range("N9").Value = [double]"1.354768"
$ws.Range("O9").Value = [double]"0.006312838586270617"
$ws.Range("P9").Value = [double]"0.006312838586270617"
$ws.Range("Q9").Value = [double]"1.036125512691556"
$ws.Range("R9").Value = [double]"9.325129614224"
$ws.Range("S9").Value = [double]"6.2151126753185E-05"
$ws.Range("T9").Value = [double]"6.2151126753185E-05"

$ws.Range("I10").Value = [double]"0.00984519497906274"
$ws.Range("J10").Value = [double]"0.00984519497906274"
$ws.Range("M10").Value = [double]"7.781650666666667"
$ws.Range("N10").Value = [double]"23.344952"
$ws.Range("O10").Value = [double]"0.1087809232135948"
$ws.Range("P10").Value = [double]"0.1087809232135948"
$ws.Range("Q10").Value = [double]"17.85420113241511"
$ws.Range("R10").Value = [double]"160.687810191736"
$ws.Range("S10").Value = [double]"0.001070969399040293"
$ws.Range("T10").Value = [double]"0.001070969399040293"

$ws.Range("I11").Value = [double]"0.00984519497906274"
$ws.Range("J11").Value = [double]"0.00984519497906274"
$ws.Range("M11").Value = [double]"3.892567333333333"
$ws.Range("N11").Value = [double]"11.677702"
$ws.Range("O11").Value = [double]"0.05441481329981927"
$ws.Range("P11").Value = [double]"0.05441481329981927"
$ws.Range("Q11").Value = [double]"8.931097406942889"
$ws.Range("R11").Value = [double]"80.379876662486"
$ws.Range("S11").Value = [double]"0.0005357244466860172"
$ws.Range("T11").Value = [double]"0.0005357244466860171"

$ws.Range("I12").Value = [double]"0.00984519497906274"
$ws.Range("J12").Value = [double]"0.00984519497906274"
$ws.Range("M12").Value = [double]"43.49559133333333"
$ws.Range("N12").Value = [double]"130.486774"
$ws.Range("O12").Value = [double]"0.6080317382054886"
$ws.Range("P12").Value = [double]"0.6080317382054886"
$ws.Range("Q12").Value = [double]"99.79618326548689"
$ws.Range("R12").Value = [double]"898.165649389382"
$ws.Range("S12").Value = [double]"0.005986191016091467"
$ws.Range("T12").Value = [double]"0.005986191016091467"

$ws.Range("I13").Value = [double]"0.00984519497906274"
$ws.Range("J13").Value = [double]"0.00984519497906274"
$ws.Range("M13").Value = [double]"7.007517"
$ws.Range("N13").Value = [double]"21.022551"
$ws.Range("O13").Value = [double]"0.09795918646930096"
$ws.Range("P13").Value = [double]"0.09795918646930096"
$ws.Range("Q13").Value = [double]"16.078030653927"
$ws.Range("R13").Value = [double]"144.702275885343"
$ws.Range("S13").Value = [double]"0.0009644272907806326"
$ws.Range("T13").Value = [double]"0.0009644272907806326"

$ws.Range("G14").Value = [double]"75.03856666666667"
$ws.Range("H14").Value = [double]"225.1157"
$ws.Range("I14").Value = [double]"0.3219883503699801"
$ws.Range("J14").Value = [double]"0.3219883503699801"
$ws.Range("M14").Value = [double]"8.906151666666668"
$ws.Range("N14").Value = [double]"26.718455"
$ws.Range("O14").Value = [double]"0.1245005002255258"
$ws.Range("P14").Value = [double]"0.1245005002255258"
$ws.Range("Q14").Value = [double]"668.3048555826113"
$ws.Range("R14").Value = [double]"6014.743700243501"
$ws.Range("S14").Value = [double]"0.04008771068785438"
$ws.Range("T14").Value = [double]"0.04008771068785438"

$ws.Range("G15").Value = [double]"75.03856666666667"
$ws.Range("H15").Value = [double]"225.1157"
$ws.Range("I15").Value = [double]"0.3219883503699801"
$ws.Range("J15").Value = [double]"0.3219883503699801"
$ws.Range("M15").Value = [double]"0.4515893333333333"
$ws.Range("N15").Value = [double]"1.354768"
$ws.Range("O15").Value = [double]"0.006312838586270617"
$ws.Range("P15").Value = [double]"0.006312838586270617"
$ws.Range("Q15").Value = [double]"33.88661629528889"
$ws.Range("R15").Value = [double]"304.9795466576"
$ws.Range("S15").Value = [double]"0.002032660482545233"
$ws.Range("T15").Value = [double]"0.002032660482545233"

$ws.Range("G16").Value = [double]"75.03856666666667"
$ws.Range("H16").Value = [double]"225.1157"
$ws.Range("I16").Value = [double]"0.3219883503699801"
$ws.Range("J16").Value = [double]"0.3219883503699801"
$ws.Range("M16").Value = [double]"7.781650666666667"
$ws.Range("N16").Value = [double]"23.344952"
$ws.Range("O16").Value = [double]"0.1087809232135948"
$ws.Range("P16").Value = [double]"0.1087809232135948"
$ws.Range("Q16").Value = [double]"583.9239123273778"
$ws.Range("R16").Value = [double]"5255.3152109464"
$ws.Range("S16").Value = [double]"0.03502619001726886"
$ws.Range("T16").Value = [double]"0.03502619001726886"

$ws.Range("G17").Value = [double]"75.03856666666667"
$ws.Range("H17").Value = [double]"225.1157"
$ws.Range("I17").Value = [double]"0.3219883503699801"
$ws.Range("J17").Value = [double]"0.3219883503699801"
$ws.Range("M17").Value = [double]"3.892567333333333"
$ws.Range("N17").Value = [double]"11.677702"
$ws.Range("O17").Value = [double]"0.05441481329981927"
$ws.Range("P17").Value = [double]"0.05441481329981927"
$ws.Range("Q17").Value = [double]"292.0926733468223"
$ws.Range("R17").Value = [double]"2628.8340601214"
$ws.Range("S17").Value = [double]"0.01752093597009926"
$ws.Range("T17").Value = [double]"0.01752093597009926"

$ws.Range("G18").Value = [double]"75.03856666666667"
$ws.Range("H18").Value = [double]"225.1157"
$ws.Range("I18").Value = [double]"0.3219883503699801"
$ws.Range("J18").Value = [double]"0.3219883503699801"
$ws.Range("M18").Value = [double]"43.49559133333333"
$ws.Range("N18").Value = [double]"130.486774"
$ws.Range("O18").Value = [double]"0.6080317382054886"
$ws.Range("P18").Value = [double]"0.6080317382054886"
$ws.Range("Q18").Value = [double]"3263.846829972422"
$ws.Range("R18").Value = [double]"29374.6214697518"
$ws.Range("S18").Value = [double]"0.1957791363573769"
$ws.Range("T18").Value = [double]"0.1957791363573769"

$ws.Range("G19").Value = [double]"75.03856666666667"
$ws.Range("H19").Value = [double]"225.1157"
$ws.Range("I19").Value = [double]"0.3219883503699801"
$ws.Range("J19").Value = [double]"0.3219883503699801"
$ws.Range("M19").Value = [double]"7.007517"
$ws.Range("N19").Value = [double]"21.022551"
$ws.Range("O19").Value = [double]"0.09795918646930096"
$ws.Range("P19").Value = [double]"0.09795918646930096"
$ws.Range("Q19").Value = [double]"525.8340315723"
$ws.Range("R19").Value = [double]"4732.5062841507"
$ws.Range("S19").Value = [double]"0.03154171685483549"
$ws.Range("T19").Value = [double]"0.03154171685483549"

$ws.Range("G20").Value = [double]"1.172599333333333"
$ws.Range("H20").Value = [double]"3.517798"
$ws.Range("I20").Value = [double]"0.005031590310914854"
$ws.Range("J20").Value = [double]"0.005031590310914854"
$ws.Range("M20").Value = [double]"8.906151666666668"
$ws.Range("N20").Value = [double]"26.718455"
$ws.Range("O20").Value = [double]"0.1245005002255258"
$ws.Range("P20").Value = [double]"0.1245005002255258"
$ws.Range("Q20").Value = [double]"10.44334750689889"
$ws.Range("R20").Value = [double]"93.99012756209001"
$ws.Range("S20").Value = [double]"0.0006264355106388082"
$ws.Range("T20").Value = [double]"0.0006264355106388081"

$ws.Range("G21").Value = [double]"1.172599333333333"
$ws.Range("H21").Value = [double]"3.517798"
$ws.Range("I21").Value = [double]"0.005031590310914854"
$ws.Range("J21").Value = [double]"0.005031590310914854"
$ws.Range("M21").Value = [double]"0.4515893333333333"
$ws.Range("N21").Value = [double]"1.354768"
$ws.Range("O21").Value = [double]"0.006312838586270617"
$ws.Range("P21").Value = [double]"0.006312838586270617"
$ws.Range("Q21").Value = [double]"0.5295333512071111"
$ws.Range("R21").Value = [double]"4.765800160864"
$ws.Range("S21").Value = [double]"3.176361746504866E-05"
$ws.Range("T21").Value = [double]"3.176361746504866E-05"

$ws.Range("G22").Value = [double]"1.172599333333333"
$ws.Range("H22").Value = [double]"3.517798"
$ws.Range("I22").Value = [double]"0.005031590310914854"
$ws.Range("J22").Value = [double]"0.005031590310914854"
$ws.Range("M22").Value = [double]"7.781650666666667"
$ws.Range("N22").Value = [double]"23.344952"
$ws.Range("O22").Value = [double]"0.1087809232135948"
$ws.Range("P22").Value = [double]"0.1087809232135948"
$ws.Range("Q22").Value = [double]"9.124758383966222"
$ws.Range("R22").Value = [double]"82.122825455696"
$ws.Range("S22").Value = [double]"0.0005473410392538964"
$ws.Range("T22").Value = [double]"0.0005473410392538964"

$ws.Range("G23").Value = [double]"1.172599333333333"
$ws.Range("H23").Value = [double]"3.517798"
$ws.Range("I23").Value = [double]"0.005031590310914854"
$ws.Range("J23").Value = [double]"0.005031590310914854"
$ws.Range("M23").Value = [double]"3.892567333333333"
$ws.Range("N23").Value = [double]"11.677702"
$ws.Range("O23").Value = [double]"0.05441481329981927"
$ws.Range("P23").Value = [double]"0.05441481329981927"
$ws.Range("Q23").Value = [double]"4.564421860021778"
$ws.Range("R23").Value = [double]"41.079796740196"
$ws.Range("S23").Value = [double]"0.0002737930473696114"
$ws.Range("T23").Value = [double]"0.0002737930473696114"

$ws.Range("G24").Value = [double]"1.172599333333333"
$ws.Range("H24").Value = [double]"3.517798"
$ws.Range("I24").Value = [double]"0.005031590310914854"
$ws.Range("J24").Value = [double]"0.005031590310914854"
$ws.Range("M24").Value = [double]"43.49559133333333"
$ws.Range("N24").Value = [double]"130.486774"
$ws.Range("O24").Value = [double]"0.6080317382054886"
$ws.Range("P24").Value = [double]"0.6080317382054886"
$ws.Range("Q24").Value = [double]"51.00290140040578"
$ws.Range("R24").Value = [double]"459.026112603652"
$ws.Range("S24").Value = [double]"0.003059366602683454"
$ws.Range("T24").Value = [double]"0.003059366602683454"

$ws.Range("G25").Value = [double]"1.172599333333333"
$ws.Range("H25").Value = [double]"3.517798"
$ws.Range("I25").Value = [double]"0.005031590310914854"
$ws.Range("J25").Value = [double]"0.005031590310914854"
$ws.Range("M25").Value = [double]"7.007517"
$ws.Range("N25").Value = [double]"21.022551"
$ws.Range("O25").Value = [double]"0.09795918646930096"
$ws.Range("P25").Value = [double]"0.09795918646930096"
$ws.Range("Q25").Value = [double]"8.217009762522"
$ws.Range("R25").Value = [double]"73.953087862698"
$ws.Range("S25").Value = [double]"0.0004928904935040362"
$ws.Range("T25").Value = [double]"0.0004928904935040362"

$ws.Range("G26").Value = [double]"115.3317996666667"
$ws.Range("H26").Value = [double]"345.995399"
$ws.Range("I26").Value = [double]"0.4948854644949822"
$ws.Range("J26").Value = [double]"0.4948854644949822"
$ws.Range("M26").Value = [double]"8.906151666666668"
$ws.Range("N26").Value = [double]"26.718455"
$ws.Range("O26").Value = [double]"0.1245005002255258"
$ws.Range("P26").Value = [double]"0.1245005002255258"
$ws.Range("Q26").Value = [double]"1027.16249982095"
$ws.Range("R26").Value = [double]"9244.462498388546"
$ws.Range("S26").Value = [double]"0.06161348788396696"
$ws.Range("T26").Value = [double]"0.06161348788396695"

$ws.Range("G27").Value = [double]"115.3317996666667"
$ws.Range("H27").Value = [double]"345.995399"
$ws.Range("I27").Value = [double]"0.4948854644949822"
$ws.Range("J27").Value = [double]"0.4948854644949822"
$ws.Range("M27").Value = [double]"0.4515893333333333"
$ws.Range("N27").Value = [double]"1.354768"
$ws.Range("O27").Value = [double]"0.006312838586270617"
$ws.Range("P27").Value = [double]"0.006312838586270617"
$ws.Range("Q27").Value = [double]"52.08261052360356"
$ws.Range("R27").Value = [double]"468.743494712432"
$ws.Range("S27").Value = [double]"0.003124132056048381"
$ws.Range("T27").Value = [double]"0.003124132056048381"

$ws.Range("G28").Value = [double]"115.3317996666667"
$ws.Range("H28").Value = [double]"345.995399"
$ws.Range("I28").Value = [double]"0.4948854644949822"
$ws.Range("J28").Value = [double]"0.4948854644949822"
$ws.Range("M28").Value = [double]"7.781650666666667"
$ws.Range("N28").Value = [double]"23.344952"
$ws.Range("O28").Value = [double]"0.1087809232135948"
$ws.Range("P28").Value = [double]"0.1087809232135948"
$ws.Range("Q28").Value = [double]"897.4717757639831"
$ws.Range("R28").Value = [double]"8077.245981875848"
$ws.Range("S28").Value = [double]"0.05383409771275285"
$ws.Range("T28").Value = [double]"0.05383409771275285"

$ws.Range("G29").Value = [double]"115.3317996666667"
$ws.Range("H29").Value = [double]"345.995399"
$ws.Range("I29").Value = [double]"0.4948854644949822"
$ws.Range("J29").Value = [double]"0.4948854644949822"
$ws.Range("M29").Value = [double]"3.892567333333333"
$ws.Range("N29").Value = [double]"11.677702"
$ws.Range("O29").Value = [double]"0.05441481329981927"
$ws.Range("P29").Value = [double]"0.05441481329981927"
$ws.Range("Q29").Value = [double]"448.9367958770109"
$ws.Range("R29").Value = [double]"4040.431162893098"
$ws.Range("S29").Value = [double]"0.02692910015528879"
$ws.Range("T29").Value = [double]"0.02692910015528879"

$ws.Range("G30").Value = [double]"115.3317996666667"
$ws.Range("H30").Value = [double]"345.995399"
$ws.Range("I30").Value = [double]"0.4948854644949822"
$ws.Range("J30").Value = [double]"0.4948854644949822"
$ws.Range("M30").Value = [double]"43.49559133333333"
$ws.Range("N30").Value = [double]"130.486774"
$ws.Range("O30").Value = [double]"0.6080317382054886"
$ws.Range("P30").Value = [double]"0.6080317382054886"
$ws.Range("Q30").Value = [double]"5016.424826039202"
$ws.Range("R30").Value = [double]"45147.82343435283"
$ws.Range("S30").Value = [double]"0.3009060691895146"
$ws.Range("T30").Value = [double]"0.3009060691895146"

$ws.Range("G31").Value = [double]"115.3317996666667"
$ws.Range("H31").Value = [double]"345.995399"
$ws.Range("I31").Value = [double]"0.4948854644949822"
$ws.Range("J31").Value = [double]"0.4948854644949822"
$ws.Range("M31").Value = [double]"7.007517"
$ws.Range("N31").Value = [double]"21.022551"
$ws.Range("O31").Value = [double]"0.09795918646930096"
$ws.Range("P31").Value = [double]"0.09795918646930096"
$ws.Range("Q31").Value = [double]"808.189546804761"
$ws.Range("R31").Value = [double]"7273.705921242849"
$ws.Range("S31").Value = [double]"0.04847857749741058"
$ws.Range("T31").Value = [double]"0.04847857749741058"

$ws.Range("G32").Value = [double]"38.39199066666666"
$ws.Range("H32").Value = [double]"115.175972"
$ws.Range("I32").Value = [double]"0.1647389374732149"
$ws.Range("J32").Value = [double]"0.1647389374732149"
$ws.Range("M32").Value = [double]"8.906151666666668"
$ws.Range("N32").Value = [double]"26.718455"
$ws.Range("O32").Value = [double]"0.1245005002255258"
$ws.Range("P32").Value = [double]"0.1245005002255258"
$ws.Range("Q32").Value = [double]"341.9248916625845"
$ws.Range("R32").Value = [double]"3077.32402496326"
$ws.Range("S32").Value = [double]"0.02051008012203687"
$ws.Range("T32").Value = [double]"0.02051008012203687"

$ws.Range("G33").Value = [double]"38.39199066666666"
$ws.Range("H33").Value = [double]"115.175972"
$ws.Range("I33").Value = [double]"0.1647389374732149"
$ws.Range("J33").Value = [double]"0.1647389374732149"
$ws.Range("M33").Value = [double]"0.4515893333333333"
$ws.Range("N33").Value = [double]"1.354768"
$ws.Range("O33").Value = [double]"0.006312838586270617"
$ws.Range("P33").Value = [double]"0.006312838586270617"
$ws.Range("Q33").Value = [double]"17.33741347049956"
$ws.Range("R33").Value = [double]"156.036721234496"
$ws.Range("S33").Value = [double]"0.001039970321142134"
$ws.Range("T33").Value = [double]"0.001039970321142134"

$ws.Range("G34").Value = [double]"38.39199066666666"
$ws.Range("H34").Value = [double]"115.175972"
$ws.Range("I34").Value = [double]"0.1647389374732149"
$ws.Range("J34").Value = [double]"0.1647389374732149"
$ws.Range("M34").Value = [double]"7.781650666666667"
$ws.Range("N34").Value = [double]"23.344952"
$ws.Range("O34").Value = [double]"0.1087809232135948"
$ws.Range("P34").Value = [double]"0.1087809232135948"
$ws.Range("Q34").Value = [double]"298.7530597659271"
$ws.Range("R34").Value = [double]"2688.777537893344"
$ws.Range("S34").Value = [double]"0.01792045370756299"
$ws.Range("T34").Value = [double]"0.01792045370756299"

$ws.Range("G35").Value = [double]"38.39199066666666"
$ws.Range("H35").Value = [double]"115.175972"
$ws.Range("I35").Value = [double]"0.1647389374732149"
$ws.Range("J35").Value = [double]"0.1647389374732149"
$ws.Range("M35").Value = [double]"3.892567333333333"
$ws.Range("N35").Value = [double]"11.677702"
$ws.Range("O35").Value = [double]"0.05441481329981927"
$ws.Range("P35").Value = [double]"0.05441481329981927"
$ws.Range("Q35").Value = [double]"149.4434087307049"
$ws.Range("R35").Value = [double]"1344.990678576344"
$ws.Range("S35").Value = [double]"0.008964238525815592"
$ws.Range("T35").Value = [double]"0.00896423852581559"

$ws.Range("G36").Value = [double]"38.39199066666666"
$ws.Range("H36").Value = [double]"115.175972"
$ws.Range("I36").Value = [double]"0.1647389374732149"
$ws.Range("J36").Value = [double]"0.1647389374732149"
$ws.Range("M36").Value = [double]"43.49559133333333"
$ws.Range("N36").Value = [double]"130.486774"
$ws.Range("O36").Value = [double]"0.6080317382054886"
$ws.Range("P36").Value = [double]"0.6080317382054886"
$ws.Range("Q36").Value = [double]"1669.882336510481"
$ws.Range("R36").Value = [double]"15028.94102859433"
$ws.Range("S36").Value = [double]"0.1001665025019642"
$ws.Range("T36").Value = [double]"0.1001665025019642"

$ws.Range("G37").Value = [double]"38.39199066666666"
$ws.Range("H37").Value = [double]"115.175972"
$ws.Range("I37").Value = [double]"0.1647389374732149"
$ws.Range("J37").Value = [double]"0.1647389374732149"
$ws.Range("M37").Value = [double]"7.007517"
$ws.Range("N37").Value = [double]"21.022551"
$ws.Range("O37").Value = [double]"0.09795918646930096"
$ws.Range("P37").Value = [double]"0.09795918646930096"
$ws.Range("Q37").Value = [double]"269.032527260508"
$ws.Range("R37").Value = [double]"2421.292745344572"
$ws.Range("S37").Value = [double]"0.01613769229469317"
$ws.Range("T37").Value = [double]"0.01613769229469317"
